$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to the new TPM-derived numbers and string labels
$ws.Range("B2").Value = "Cort"
$ws.Range("C2").Value = "Sstr1"
$ws.Range("D2").Value = "ECs"

$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.1365053333333333
$ws.Range("N2").Value = 0.409516
$ws.Range("Q2").Value = 0.03127387238622222
$ws.Range("R2").Value = 0.281464851476
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove rows 3 and 4 entirely (data now only has a single data row)
$ws.Rows("3:4").Delete()
